$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.747.22"
$ws.Range("E2").Value = "  +1.78%  "
$ws.Range("D3").Value = "1.655.13"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.03"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3823"
$ws.Range("E7").Value = "  +1.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3612"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.18"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.255"
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08216"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.72"
$ws.Range("E13").Value = "  +2.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.549"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.453"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001240"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "1.637.92"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.81"
$ws.Range("E18").Value = "  +3.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06984"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.789"
$ws.Range("E20").Value = "  +3.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.78"
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.75"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").Value = "23.752.89"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.565"
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.087"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.92"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.224"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").Value = "1.820.79"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.924"
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.167"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.083"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.90"
$ws.Range("E35").Value = "  +5.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02835"
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.159"
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2521"
$ws.Range("E38").Value = "  +1.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08832"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.07183"
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "13.08"
$ws.Range("E41").Value = "  +8.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7073"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.340"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.95"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6549"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.332"
$ws.Range("E46").Value = "  +3.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.960"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07976"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.61"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("E51").Value = "  +1.03%  "
